$d = $word.ActiveDocument

$pairs = @(
    @("/vremenskog kontingenta. Ovaj film ", "/vremenskog kontingenta. ## Ovaj film "),
    @(" Einsteinovoj teoriji opće relativnosti. Film je metodički režiran, ", " Einsteinovoj teoriji opće relativnosti. ## Film je metodički režiran, "),
    @(" detalje čovjeka (Hawking) kao i njegov rad (Crne rupe). Intervjui sa svojom obitelji su malo ", " detalje čovjeka (Hawking) kao i njegov rad (Crne rupe). ## Intervjui sa svojom obitelji su malo "),
    @("njegove teorije i ideje. Philip Glass soundtrack ", "njegove teorije i ideje. ## Philip Glass soundtrack "),
    @(" pohvale filmu. Samo je jedan drugi čovjek mogao skladati takve uklete instelarne melodije (Jean Michel Jarre). Sveukupno bih visoko preporučio ovaj film na temelju Hawkingovih '", " pohvale filmu. ## Samo je jedan drugi čovjek mogao skladati takve uklete instelarne melodije (Jean Michel Jarre). ## Sveukupno bih visoko preporučio ovaj film na temelju Hawkingovih '"),
    @("tako grozan film u dugo... dugo vremena... ", "tako grozan film u dugo... dugo vremena...  ## "),
    @("sam ga sinoć i htio otići nakon 20 minuta... Keira Knightley pokušava stvarno ", "sam ga sinoć i htio otići nakon 20 minuta... ## Keira Knightley pokušava stvarno "),
    @("u neko vrijeme i nije imao karizmu ispuniti ulogu... Sienna Millers ", "u neko vrijeme i nije imao karizmu ispuniti ulogu... ## Sienna Millers "),
    @(": Je li ikad imala satove glume? Sudeći po ", ": Je li ikad imala satove glume? ## Sudeći po "),
    @("u bliskoj budućnosti... oboje izgledaju jako lijepo.. možda je to ono ", "u bliskoj budućnosti... ## oboje izgledaju jako lijepo.. možda je to ono "),
    @("buduću karijeru.. ako ", "buduću karijeru.. ## ako "),
    @("Zahvaljujući drugim recenzentima koji su me usmjerili na ovaj proizvod kad mi je rečeno da sam anemična. Sada uzimam ove ", "Zahvaljujući drugim recenzentima koji su me usmjerili na ovaj proizvod kad mi je rečeno da sam anemična. ## Sada uzimam ove "),
    @(" oko 4 mjeseca, a anemija je nestala. Dobar proizvod. ", " oko 4 mjeseca, a anemija je nestala. ## Dobar proizvod. ## "),
    @("Ovo je jedan od mojih omiljenih deserta i brzo se topi u ustima. Ova marka je dobra i ", "Ovo je jedan od mojih omiljenih deserta i brzo se topi u ustima. ## Ova marka je dobra i "),
    @(" dobro zapakirana. Svatko bi trebao probati ", " dobro zapakirana. ## Svatko bi trebao probati "),
    @(" jednom. Cijena ", " jednom. ## Cijena "),
    @("/poklon za mlade i stare. ", "/poklon za mlade i stare. ## "),
    @(" broj načina. Sjajno je i teško ćeš ga držati podalje od odraslih.", " broj načina. ## Sjajno je i teško ćeš ga držati podalje od odraslih."),
    @("Još jedan loš zombi film. U usporedbi s većinom drugih, jedina razlika ovdje je ", "Još jedan loš zombi film. ## U usporedbi s većinom drugih, jedina razlika ovdje je "),
    @(" žensko. ", " žensko. ## "),
    @("je ista. Akcijske scene nisu zanimljive. Specijalni efekti su ", "je ista. ## Akcijske scene nisu zanimljive. ## Specijalni efekti su "),
    @("Yoyo izgleda izvan ravnoteže. Bez obzira na ", "Yoyo izgleda izvan ravnoteže. ## Bez obzira na "),
    @(" na jednu stranu. ", " na jednu stranu. ## "),
    @(". Imam malo iskustva s ", ". ## Imam malo iskustva s "),
    @(". Nemojte ih ", ". ## Nemojte ih "),
    @(". Oni jednostavno neće štititi/podupirati vaše zglobove ", ". ## Oni jednostavno neće štititi/podupirati vaše zglobove ")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Text = $new
    } else {
        Write-Output "NOT FOUND: $old"
    }
}
